$d = $word.ActiveDocument

# Locate the paragraph "One manuscript in-review." (currently split across
# four runs: "One manuscript " / "in" / "-" / "review.") and collapse a
# Range onto exactly that text span.
$target = $d.Content
$found = $target.Find.Execute("One manuscript in-review.", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'One manuscript in-review.' paragraph text"
}

# Re-anchor onto a fresh Range over the exact same span - InsertXML needs to be
# invoked on a plain start/end Range (not the live Find range) to replace the
# found text rather than insert alongside it.
$span = $d.Range($target.Start, $target.End)

# Build the replacement content as its own distinct runs - "Two" / " manuscript" /
# "s" / " " / "accepted for Conference" / "." - by inserting a WordprocessingML
# fragment over the located range, which preserves run boundaries (unlike
# Range.Text / Range.InsertAfter, which coalesce adjacently-inserted text that
# shares formatting into a single run).
$paragraphXml = '<w:p>' + `
    '<w:r><w:t>Two</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> manuscript</w:t></w:r>' + `
    '<w:r><w:t>s</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>accepted for Conference</w:t></w:r>' + `
    '<w:r><w:t>.</w:t></w:r>' + `
    '</w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
            '<w:body>' + $paragraphXml + '</w:body>' + `
          '</w:document>' + `
        '</pkg:xmlData>' + `
      '</pkg:part>' + `
    '</pkg:package>'

$span.InsertXML($packageXml)
